$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(120, 1).Formula = "'2026-01-28"
$ws.Cells.Item(120, 2).Value = "15:38:21"
$ws.Cells.Item(120, 3).Value = "15:00"
$ws.Cells.Item(120, 4).Value = "Bathroom"
$ws.Cells.Item(120, 5).Value = "No Motion"
$ws.Cells.Item(120, 6).Value = "Inactive"

$ws.Cells.Item(121, 1).Formula = "'2026-01-28"
$ws.Cells.Item(121, 2).Value = "15:38:22"
$ws.Cells.Item(121, 3).Value = "15:00"
$ws.Cells.Item(121, 4).Value = "Bathroom"
$ws.Cells.Item(121, 5).Value = "No Motion"
$ws.Cells.Item(121, 6).Value = "Inactive"

$ws.Cells.Item(122, 1).Formula = "'2026-01-28"
$ws.Cells.Item(122, 2).Value = "15:38:26"
$ws.Cells.Item(122, 3).Value = "15:00"
$ws.Cells.Item(122, 4).Value = "Bathroom"
$ws.Cells.Item(122, 5).Value = "No Motion"
$ws.Cells.Item(122, 6).Value = "Inactive"

$ws.Cells.Item(123, 1).Formula = "'2026-01-28"
$ws.Cells.Item(123, 2).Value = "15:38:31"
$ws.Cells.Item(123, 3).Value = "15:00"
$ws.Cells.Item(123, 4).Value = "Bathroom"
$ws.Cells.Item(123, 5).Value = "No Motion"
$ws.Cells.Item(123, 6).Value = "Inactive"

$ws.Cells.Item(124, 1).Formula = "'2026-01-28"
$ws.Cells.Item(124, 2).Value = "15:38:36"
$ws.Cells.Item(124, 3).Value = "15:00"
$ws.Cells.Item(124, 4).Value = "Bathroom"
$ws.Cells.Item(124, 5).Value = "No Motion"
$ws.Cells.Item(124, 6).Value = "Inactive"

$ws.Cells.Item(125, 1).Formula = "'2026-01-28"
$ws.Cells.Item(125, 2).Value = "15:38:41"
$ws.Cells.Item(125, 3).Value = "15:00"
$ws.Cells.Item(125, 4).Value = "Bathroom"
$ws.Cells.Item(125, 5).Value = "No Motion"
$ws.Cells.Item(125, 6).Value = "Inactive"

$ws.Cells.Item(126, 1).Formula = "'2026-01-28"
$ws.Cells.Item(126, 2).Value = "15:38:46"
$ws.Cells.Item(126, 3).Value = "15:00"
$ws.Cells.Item(126, 4).Value = "Bathroom"
$ws.Cells.Item(126, 5).Value = "No Motion"
$ws.Cells.Item(126, 6).Value = "Inactive"

$ws.Cells.Item(127, 1).Formula = "'2026-01-28"
$ws.Cells.Item(127, 2).Value = "15:38:51"
$ws.Cells.Item(127, 3).Value = "15:00"
$ws.Cells.Item(127, 4).Value = "Bathroom"
$ws.Cells.Item(127, 5).Value = "No Motion"
$ws.Cells.Item(127, 6).Value = "Inactive"

$ws.Cells.Item(128, 1).Formula = "'2026-01-28"
$ws.Cells.Item(128, 2).Value = "15:38:56"
$ws.Cells.Item(128, 3).Value = "15:00"
$ws.Cells.Item(128, 4).Value = "Bathroom"
$ws.Cells.Item(128, 5).Value = "No Motion"
$ws.Cells.Item(128, 6).Value = "Inactive"

$ws.Cells.Item(129, 1).Formula = "'2026-01-28"
$ws.Cells.Item(129, 2).Value = "15:39:01"
$ws.Cells.Item(129, 3).Value = "15:00"
$ws.Cells.Item(129, 4).Value = "Bathroom"
$ws.Cells.Item(129, 5).Value = "No Motion"
$ws.Cells.Item(129, 6).Value = "Inactive"

$ws.Cells.Item(130, 1).Formula = "'2026-01-28"
$ws.Cells.Item(130, 2).Value = "15:39:06"
$ws.Cells.Item(130, 3).Value = "15:00"
$ws.Cells.Item(130, 4).Value = "Bathroom"
$ws.Cells.Item(130, 5).Value = "No Motion"
$ws.Cells.Item(130, 6).Value = "Inactive"

$ws.Cells.Item(131, 1).Formula = "'2026-01-28"
$ws.Cells.Item(131, 2).Value = "15:39:11"
$ws.Cells.Item(131, 3).Value = "15:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "No Motion"
$ws.Cells.Item(131, 6).Value = "Inactive"

$ws.Cells.Item(132, 1).Formula = "'2026-01-28"
$ws.Cells.Item(132, 2).Value = "15:39:17"
$ws.Cells.Item(132, 3).Value = "15:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "No Motion"
$ws.Cells.Item(132, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(121, 1).Formula = "'2026-01-28"
$ws.Cells.Item(121, 2).Value = "15:38:21"
$ws.Cells.Item(121, 3).Value = "15:00"
$ws.Cells.Item(121, 4).Value = "Bathroom"
$ws.Cells.Item(121, 5).Formula = "'87.5%"
$ws.Cells.Item(121, 6).Value = "Active"

$ws.Cells.Item(122, 1).Formula = "'2026-01-28"
$ws.Cells.Item(122, 2).Value = "15:38:24"
$ws.Cells.Item(122, 3).Value = "15:00"
$ws.Cells.Item(122, 4).Value = "Bathroom"
$ws.Cells.Item(122, 5).Formula = "'88.4%"
$ws.Cells.Item(122, 6).Value = "Active"

$ws.Cells.Item(123, 1).Formula = "'2026-01-28"
$ws.Cells.Item(123, 2).Value = "15:38:32"
$ws.Cells.Item(123, 3).Value = "15:00"
$ws.Cells.Item(123, 4).Value = "Bathroom"
$ws.Cells.Item(123, 5).Formula = "'87.5%"
$ws.Cells.Item(123, 6).Value = "Active"

$ws.Cells.Item(124, 1).Formula = "'2026-01-28"
$ws.Cells.Item(124, 2).Value = "15:38:36"
$ws.Cells.Item(124, 3).Value = "15:00"
$ws.Cells.Item(124, 4).Value = "Bathroom"
$ws.Cells.Item(124, 5).Formula = "'88.5%"
$ws.Cells.Item(124, 6).Value = "Active"

$ws.Cells.Item(125, 1).Formula = "'2026-01-28"
$ws.Cells.Item(125, 2).Value = "15:38:44"
$ws.Cells.Item(125, 3).Value = "15:00"
$ws.Cells.Item(125, 4).Value = "Bathroom"
$ws.Cells.Item(125, 5).Formula = "'87.0%"
$ws.Cells.Item(125, 6).Value = "Active"

$ws.Cells.Item(126, 1).Formula = "'2026-01-28"
$ws.Cells.Item(126, 2).Value = "15:38:52"
$ws.Cells.Item(126, 3).Value = "15:00"
$ws.Cells.Item(126, 4).Value = "Bathroom"
$ws.Cells.Item(126, 5).Formula = "'88.5%"
$ws.Cells.Item(126, 6).Value = "Active"

$ws.Cells.Item(127, 1).Formula = "'2026-01-28"
$ws.Cells.Item(127, 2).Value = "15:38:56"
$ws.Cells.Item(127, 3).Value = "15:00"
$ws.Cells.Item(127, 4).Value = "Bathroom"
$ws.Cells.Item(127, 5).Formula = "'88.5%"
$ws.Cells.Item(127, 6).Value = "Active"

$ws.Cells.Item(128, 1).Formula = "'2026-01-28"
$ws.Cells.Item(128, 2).Value = "15:39:00"
$ws.Cells.Item(128, 3).Value = "15:00"
$ws.Cells.Item(128, 4).Value = "Bathroom"
$ws.Cells.Item(128, 5).Formula = "'87.6%"
$ws.Cells.Item(128, 6).Value = "Active"

$ws.Cells.Item(129, 1).Formula = "'2026-01-28"
$ws.Cells.Item(129, 2).Value = "15:39:04"
$ws.Cells.Item(129, 3).Value = "15:00"
$ws.Cells.Item(129, 4).Value = "Bathroom"
$ws.Cells.Item(129, 5).Formula = "'88.5%"
$ws.Cells.Item(129, 6).Value = "Active"

$ws.Cells.Item(130, 1).Formula = "'2026-01-28"
$ws.Cells.Item(130, 2).Value = "15:39:08"
$ws.Cells.Item(130, 3).Value = "15:00"
$ws.Cells.Item(130, 4).Value = "Bathroom"
$ws.Cells.Item(130, 5).Formula = "'88.5%"
$ws.Cells.Item(130, 6).Value = "Active"

$ws.Cells.Item(131, 1).Formula = "'2026-01-28"
$ws.Cells.Item(131, 2).Value = "15:39:12"
$ws.Cells.Item(131, 3).Value = "15:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Formula = "'87.5%"
$ws.Cells.Item(131, 6).Value = "Active"

$ws.Cells.Item(132, 1).Formula = "'2026-01-28"
$ws.Cells.Item(132, 2).Value = "15:39:16"
$ws.Cells.Item(132, 3).Value = "15:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Formula = "'88.4%"
$ws.Cells.Item(132, 6).Value = "Active"

$ws.Cells.Item(133, 1).Formula = "'2026-01-28"
$ws.Cells.Item(133, 2).Value = "15:39:20"
$ws.Cells.Item(133, 3).Value = "15:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Formula = "'87.5%"
$ws.Cells.Item(133, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(121, 1).Formula = "'2026-01-28"
$ws.Cells.Item(121, 2).Value = "15:38:22"
$ws.Cells.Item(121, 3).Value = "15:00"
$ws.Cells.Item(121, 4).Value = "Bathroom"
$ws.Cells.Item(121, 5).Value = "22.9C"
$ws.Cells.Item(121, 6).Value = "Active"

$ws.Cells.Item(122, 1).Formula = "'2026-01-28"
$ws.Cells.Item(122, 2).Value = "15:38:24"
$ws.Cells.Item(122, 3).Value = "15:00"
$ws.Cells.Item(122, 4).Value = "Bathroom"
$ws.Cells.Item(122, 5).Value = "22.9C"
$ws.Cells.Item(122, 6).Value = "Active"

$ws.Cells.Item(123, 1).Formula = "'2026-01-28"
$ws.Cells.Item(123, 2).Value = "15:38:32"
$ws.Cells.Item(123, 3).Value = "15:00"
$ws.Cells.Item(123, 4).Value = "Bathroom"
$ws.Cells.Item(123, 5).Value = "22.9C"
$ws.Cells.Item(123, 6).Value = "Active"

$ws.Cells.Item(124, 1).Formula = "'2026-01-28"
$ws.Cells.Item(124, 2).Value = "15:38:36"
$ws.Cells.Item(124, 3).Value = "15:00"
$ws.Cells.Item(124, 4).Value = "Bathroom"
$ws.Cells.Item(124, 5).Value = "22.9C"
$ws.Cells.Item(124, 6).Value = "Active"

$ws.Cells.Item(125, 1).Formula = "'2026-01-28"
$ws.Cells.Item(125, 2).Value = "15:38:44"
$ws.Cells.Item(125, 3).Value = "15:00"
$ws.Cells.Item(125, 4).Value = "Bathroom"
$ws.Cells.Item(125, 5).Value = "22.9C"
$ws.Cells.Item(125, 6).Value = "Active"

$ws.Cells.Item(126, 1).Formula = "'2026-01-28"
$ws.Cells.Item(126, 2).Value = "15:38:52"
$ws.Cells.Item(126, 3).Value = "15:00"
$ws.Cells.Item(126, 4).Value = "Bathroom"
$ws.Cells.Item(126, 5).Value = "22.9C"
$ws.Cells.Item(126, 6).Value = "Active"

$ws.Cells.Item(127, 1).Formula = "'2026-01-28"
$ws.Cells.Item(127, 2).Value = "15:38:56"
$ws.Cells.Item(127, 3).Value = "15:00"
$ws.Cells.Item(127, 4).Value = "Bathroom"
$ws.Cells.Item(127, 5).Value = "22.9C"
$ws.Cells.Item(127, 6).Value = "Active"

$ws.Cells.Item(128, 1).Formula = "'2026-01-28"
$ws.Cells.Item(128, 2).Value = "15:39:00"
$ws.Cells.Item(128, 3).Value = "15:00"
$ws.Cells.Item(128, 4).Value = "Bathroom"
$ws.Cells.Item(128, 5).Value = "22.9C"
$ws.Cells.Item(128, 6).Value = "Active"

$ws.Cells.Item(129, 1).Formula = "'2026-01-28"
$ws.Cells.Item(129, 2).Value = "15:39:04"
$ws.Cells.Item(129, 3).Value = "15:00"
$ws.Cells.Item(129, 4).Value = "Bathroom"
$ws.Cells.Item(129, 5).Value = "22.9C"
$ws.Cells.Item(129, 6).Value = "Active"

$ws.Cells.Item(130, 1).Formula = "'2026-01-28"
$ws.Cells.Item(130, 2).Value = "15:39:08"
$ws.Cells.Item(130, 3).Value = "15:00"
$ws.Cells.Item(130, 4).Value = "Bathroom"
$ws.Cells.Item(130, 5).Value = "22.9C"
$ws.Cells.Item(130, 6).Value = "Active"

$ws.Cells.Item(131, 1).Formula = "'2026-01-28"
$ws.Cells.Item(131, 2).Value = "15:39:12"
$ws.Cells.Item(131, 3).Value = "15:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "22.9C"
$ws.Cells.Item(131, 6).Value = "Active"

$ws.Cells.Item(132, 1).Formula = "'2026-01-28"
$ws.Cells.Item(132, 2).Value = "15:39:16"
$ws.Cells.Item(132, 3).Value = "15:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "22.9C"
$ws.Cells.Item(132, 6).Value = "Active"

$ws.Cells.Item(133, 1).Formula = "'2026-01-28"
$ws.Cells.Item(133, 2).Value = "15:39:20"
$ws.Cells.Item(133, 3).Value = "15:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "22.9C"
$ws.Cells.Item(133, 6).Value = "Active"
